$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.933.21"
$ws.Range("E2").Value = "  -2.04%  "

# Row 3
$ws.Range("D3").Value = "1.867.32"
$ws.Range("E3").Value = "  -2.43%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'311.99"
$ws.Range("E5").Value = "  -1.02%  "

# Row 6
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.03%  "

# Row 7
$ws.Range("D7").Value = "'0.4978"
$ws.Range("E7").Value = "  -3.58%  "

# Row 8
$ws.Range("D8").Value = "'0.3801"
$ws.Range("E8").Value = "  -4.38%  "

# Row 9
$ws.Range("D9").Value = "'0.08943"
$ws.Range("E9").Value = "  -9.36%  "

# Row 10
$ws.Range("E10").Value = "  -2.90%  "

# Row 11
$ws.Range("D11").Value = "'41.48"
$ws.Range("E11").Value = "  -1.87%  "

# Row 12
$ws.Range("D12").Value = "'6.304"
$ws.Range("E12").Value = "  -3.45%  "

# Row 13
$ws.Range("D13").Value = "'20.64"
$ws.Range("E13").Value = "  -2.45%  "

# Row 14
$ws.Range("D14").Value = "1.860.20"
$ws.Range("E14").Value = "  -3.19%  "

# Row 15
$ws.Range("D15").Value = "'7.220"
$ws.Range("E15").Value = "  -3.30%  "

# Row 16
$ws.Range("E16").Value = "  +0.01%  "

# Row 17
$ws.Range("D17").Value = "'0.00001099"
$ws.Range("E17").Value = "  -3.47%  "

# Row 18
$ws.Range("D18").Value = "'90.75"
$ws.Range("E18").Value = "  -4.05%  "

# Row 19
$ws.Range("D19").Value = "'0.06612"
$ws.Range("E19").Value = "  -0.69%  "

# Row 20
$ws.Range("D20").Value = "'17.87"
$ws.Range("E20").Value = "  -2.01%  "

# Row 21
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.00%  "

# Row 22
$ws.Range("D22").Value = "'6.087"
$ws.Range("E22").Value = "  -3.58%  "

# Row 23
$ws.Range("D23").Value = "27.963.50"
$ws.Range("E23").Value = "  -2.13%  "

# Row 24
$ws.Range("D24").Value = "'11.40"
$ws.Range("E24").Value = "  -0.74%  "

# Row 25
$ws.Range("D25").Value = "'2.286"
$ws.Range("E25").Value = "  -1.60%  "

# Row 26
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "'3.384"
$ws.Range("E26").Value = "  +0.08%  "

# Row 27
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.074.66"
$ws.Range("E27").Value = "  -2.65%  "

# Row 28
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.523"
$ws.Range("E28").Value = "  -5.93%  "

# Row 29
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'157.71"
$ws.Range("E29").Value = "  +0.06%  "

# Row 30
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'20.67"
$ws.Range("E30").Value = "  -2.82%  "

# Row 31
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "'126.13"
$ws.Range("E31").Value = "  -2.43%  "

# Row 32
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.1055"
$ws.Range("E32").Value = "  -1.95%  "

# Row 33
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'1.056"
$ws.Range("E33").Value = "  -5.03%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.574"
$ws.Range("E34").Value = "  -2.92%  "

# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'3.586"
$ws.Range("E35").Value = "  -1.24%  "

# Row 36
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "'9.336"
$ws.Range("E36").Value = "  -5.60%  "

# Row 37
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06534"
$ws.Range("E37").Value = "  -3.38%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02402"
$ws.Range("E38").Value = "  -1.49%  "

# Row 39
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.2183"
$ws.Range("E39").Value = "  -1.59%  "

# Row 40
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.276"
$ws.Range("E40").Value = "  +7.47%  "

# Row 41
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'1.199"
$ws.Range("E41").Value = "  -5.76%  "

# Row 42
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'11.63"
$ws.Range("E42").Value = "  -1.37%  "

# Row 43
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.6356"
$ws.Range("E43").Value = "  -1.91%  "

# Row 44
$ws.Range("B44").Value = "InternetComputer(DFINITY)"
$ws.Range("C44").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D44").Value = "'4.893"
$ws.Range("E44").Value = "  -3.93%  "

# Row 45
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "'1.0000"
$ws.Range("E45").Value = "  -0.01%  "

# Row 46
$ws.Range("D46").Value = "'13.21"
$ws.Range("E46").Value = "  -2.91%  "

# Row 47
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5984"
$ws.Range("E47").Value = "  -1.91%  "

# Row 48
$ws.Range("B48").Value = "WEMIXTOKEN"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "'1.283"
$ws.Range("E48").Value = "  -0.47%  "

# Row 49
$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D49").Value = "'3.667"
$ws.Range("E49").Value = "  -2.66%  "

# Row 50
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "'1.222"
$ws.Range("E50").Value = "  +1.46%  "

# Row 51
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.965"
$ws.Range("E51").Value = "  -4.18%  "
